$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Al Prefetto di <$COMUNE_COMANDO>" table row (wrong
#    recipient that is being removed from the letterhead table).
# ------------------------------------------------------------------
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($ri = $tbl.Rows.Count; $ri -ge 1; $ri--) {
        $rowText = $tbl.Rows.Item($ri).Range.Text
        if ($rowText -like "*Prefetto*") {
            $tbl.Rows.Item($ri).Delete()
        }
    }
}

# ------------------------------------------------------------------
# 2) Add a new paragraph, right after the "Avverso al presente
#    parere..." paragraph, informing that the communication is also
#    sent to the Sindaco pursuant to art. 13 DPR 577/82.
# ------------------------------------------------------------------
$egrave = [char]0x00E8
$rsquo  = [char]0x2019
$newText = "Questa comunicazione " + $egrave + " trasmessa al Sindaco ai sensi dell" + $rsquo + "art. 13 del DPR 577/82."

$paras = $d.Content.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Avverso al presente parere*") {
        $p.Range.InsertParagraphAfter()
        $paras2 = $d.Content.Paragraphs
        $newPara = $paras2.Item($i + 1)
        $newRange = $newPara.Range.Duplicate

        $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
               "<w:pPr><w:pStyle w:val='Corpodeltesto'/><w:rPr><w:lang w:val='it-IT'/></w:rPr></w:pPr>" +
               "<w:r><w:rPr><w:lang w:val='it-IT'/></w:rPr><w:t xml:space='preserve'>" + $newText + "</w:t></w:r>" +
               "</w:p>"

        $null = $newRange.InsertXML($xml)
        break
    }
}
